$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the split points inside the single run-of-runs paragraph.
#    Paragraph 1 currently holds 8 runs concatenated back-to-back; the
#    target splits that text into three paragraphs (each separated by a
#    blank paragraph):
#      P1 = runs 0-2  (intro / Polaris / CBS)
#      P2 = runs 3-6  (MarketWatch / Post / awards / lecturer)
#      P3 = run 7     (Harvard degree)
# ---------------------------------------------------------------------
$full = $d.Content.Text
$splitA = $full.IndexOf('在加入 CBS 之前的十年里')
$splitB = $full.IndexOf('Larry 拥有哈佛大学工商管理硕士学位')

# Insert the paragraph breaks (back-to-front so earlier offsets remain valid).
$rB = $d.Range($splitB, $splitB)
$rB.InsertParagraphBefore()

$rA = $d.Range($splitA, $splitA)
$rA.InsertParagraphBefore()

# ---------------------------------------------------------------------
# 2. Merge the runs inside each of the three paragraphs into a single
#    run by rewriting each paragraph's text back onto itself via XML
#    (collapses the run-series into one <w:r>) and swap the paragraph
#    spacing/format for a first-line indent.
# ---------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p1 = $d.Paragraphs(1)
$t1 = $p1.Range.Text
$t1 = $t1.Substring(0, $t1.Length - 1)
$frag1 = '<w:p ' + $wNs + '><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">' + $t1 + '</w:t></w:r></w:p>'
$null = $p1.Range.InsertXML($frag1)

$p2 = $d.Paragraphs(2)
$t2 = $p2.Range.Text
$t2 = $t2.Substring(0, $t2.Length - 1)
$frag2 = '<w:p ' + $wNs + '><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">' + $t2 + '</w:t></w:r></w:p>'
$null = $p2.Range.InsertXML($frag2)

$p3 = $d.Paragraphs(3)
$t3 = $p3.Range.Text
$t3 = $t3.Substring(0, $t3.Length - 1)
$frag3 = '<w:p ' + $wNs + '><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">' + $t3 + '</w:t></w:r></w:p>'
$null = $p3.Range.InsertXML($frag3)

# ---------------------------------------------------------------------
# 3. Insert the blank separator paragraphs between 1/2 and 2/3.
# ---------------------------------------------------------------------
$d.Paragraphs(1).Range.InsertParagraphAfter()
$d.Paragraphs(3).Range.InsertParagraphAfter()

Write-Output $d.Paragraphs.Count

# ---------------------------------------------------------------------
# 4. Add the _GoBack bookmark at the very start of the final paragraph.
# ---------------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$bmRange = $d.Range($last.Range.Start, $last.Range.Start)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
